# "Updated symbol list" run — refreshes scraped coinranking.com prices
# (column D), a couple of volume/"Worst in 24h" labels (column E), and
# re-sorts three rows (41-43) whose coins swapped rank order, each with a
# freshly scraped price.
#
# All of these cells hold their numbers as literal TEXT (the sheet has no
# number formatting on column D — every price was scraped straight into a
# string cell), so plain numeric-looking assignments are forced through
# NumberFormat "@" (Text) first and the style is reset back to Normal
# afterwards so no stray formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-Label($addr, $val) {
    $ws.Range($addr).Value = $val
}

# --- Price (column D) refreshes for existing rows ---------------------
Set-TextValue "D2"  "248.06"
Set-TextValue "D4"  "5.373"
Set-TextValue "D5"  "0.05702"
Set-TextValue "D6"  "3.412"
Set-TextValue "D7"  "6.314"
Set-TextValue "D8"  "0.8111"
Set-TextValue "D9"  "0.9139"
Set-TextValue "D10" "0.1414"
Set-TextValue "D11" "0.07423"
Set-TextValue "D13" "0.03021"
Set-TextValue "D14" "0.09350"
Set-TextValue "D15" "3.718"
Set-TextValue "D16" "0.001586"
Set-TextValue "D17" "0.04769"
Set-TextValue "D18" "0.01830"

# --- Row 19 (One / ONE) picks up a "Worst in 24h" tag ------------------
Set-Label "E19" "18OneONEWorstin24h"

Set-TextValue "D20" "0.006500"
Set-TextValue "D21" "0.004999"
Set-TextValue "D22" "0.001027"
Set-TextValue "D23" "0.0001502"
Set-TextValue "D24" "3.696"
Set-TextValue "D40" "0.03984"

# --- Rows 41-43 re-sorted: KickToken, BKEXToken, CEJI (in that order) -
Set-Label "B41" "KickToken"
Set-Label "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006877"
Set-Label "E41" "40KickTokenKICK"

Set-Label "B42" "BKEXToken"
Set-Label "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1063"
Set-Label "E42" "41BKEXTokenBKK"

Set-Label "B43" "CEJI"
Set-Label "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002713"
Set-Label "E43" "42CEJICEJI"

Set-TextValue "D44" "0.007455"
